$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = "questions = [`n    {`n        `"title`": `"What is \u2018bug severity\u2019?`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"The degree of impact a defect has on the system.`",`n            `"The order in which a defect should be fixed.`",`n            `"The likelihood and the impact of the bug.`",`n            `"The cost of the impact combined with the cost to resolve the consequences of the impact.`"`n        ],`n        `"score`": `"The degree of impact a defect has on the system.`"`n    },`n    {`n        `"title`": `"A walkthrough is led by a trained moderator, whereas an inspection is led by the author.`",`n        `"ques_type`": 11,`n        `"options`": [`n            `"true`",`n            `"false`"`n        ],`n        `"score`": `"False`"`n    },`n    {`n        `"title`": `"What is a typical objective of an unit test?`",`n        `"ques_type`": 2,`n        `"options`": [`n            `"To validate that each unit of the software performs as designed.`",`n            `"To test all branches of the software program once.`",`n            `"To build confidence in the quality of interfaces.`",`n            `"To verify that the conditions match the product requirements.`"`n        ],`n        `"score`": `"To validate that each unit of the software performs as designed.`"`n    },`n    {`n        `"title`": `"Statement 1:Regression testing ensures the recent changes in the software code have not affected the unaltered sections. Retesting establishes that failed test cases are passed after the defect is fixed.Statement 2:Regression testing includes no defect verification, retesting does.Are the statements true or false?`",`n        `"ques_type`": 11,`n        `"options`": [`n            `"true`",`n            `"false`"`n        ],`n        `"score`": `"True`"`n    }`n]"

$ws.Range("A2").ClearContents()
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $text
$ws.Rows(1).AutoFit()
